$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove hyperlink before shifting rows (ref is still accurate at this point)
$ws.Range("G2").Hyperlinks.Delete()

# Insert a new row above row 2 -> becomes the sub-header row for the merged Name column
$ws.Range("A2").EntireRow.Insert()

# New sub-header row (row 2) values
$ws.Range("C2").Value = "First 1"
$ws.Range("D2").Value = "Last"

# Update header row: merge C1:D1 into "Name" (replacing "First Name"/"Last Name")
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = ""
$ws.Range("C1:D1").Merge()

# Fix the data further down: K5 (was row4 K, "invalid") becomes "неверно"
$ws.Range("K5").Value = "неверно"

# Re-add hyperlink at new location (row shifted from 2 -> 3)
$ws.Hyperlinks.Add($ws.Range("G3"), "mailto:john@examle.com")

Write-Host "done"
